# Update the "daftar_obat" query/answer table on Sheet1: several rows get
# their query text and/or drug-list text expanded with extra detail
# (e.g. additional drug names, more descriptive symptom wording).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C3').Value = 'Paracetamol, Novagesic, topcilin'
$ws.Range('B5').Value = 'rasa nyeri haid di hari pertama'
$ws.Range('C5').Value = 'Lapistan, Spasminal, Freedol'
$ws.Range('B6').Value = 'sesak nafas karena serangan asma'
$ws.Range('C6').Value = 'Aminophylline, Salbutamol Tablet, Bufakris syrup'
$ws.Range('B7').Value = 'BAB encer atau diare akut'
$ws.Range('C7').Value = 'Loperamide, Attalpugite, Bidium'
$ws.Range('B8').Value = 'gatal seperti gejala alergi'
$ws.Range('B9').Value = 'mual - mual seperti gejala asam lambung'
$ws.Range('C9').Value = 'Novamag, Berlosid, Samtacid'
$ws.Range('B10').Value = 'kekurangan darah atau anemia dan kekurangan zat besi'
$ws.Range('C11').Value = 'Ketoconazole, Miconazole, Mycostop'
$ws.Range('C12').Value = 'Vitamin B Complex, Selkom-c, Benovit C, Imnusive for Adults'
$ws.Range('C13').Value = 'Allopurinol, Omeric, Alodan'
$ws.Range('C14').Value = 'Simvastatin, Atorvastatin, Konilife redaxin'
$ws.Range('C15').Value = 'Furosemide, Farmoten, Prix'
$ws.Range('B16').Value = 'sering kesemutan akibat kekurangan vitamin b'
$ws.Range('C17').Value = 'Kurkumex, Damuvit'
$ws.Range('B18').Value = 'sakit kepala atau pusing sampai berputar seperti vertigo'
$ws.Range('C18').Value = 'Histigo, Dimenhydrinate, Merislon'
$ws.Rows.Item(18).RowHeight = 31.2
$ws.Range('B19').Value = 'saya mual dan muntah - muntah'
$ws.Range('B20').Value = 'Saya merasakan sakit kepala dan demam'
$ws.Range('C20').Value = 'topcilin, Paracetamol, Novagesic'
$ws.Range('B23').Value = 'saya nyeri ringan pada badan'
$ws.Range('C23').Value = 'Freedol, Natrium Diklofenak Tablet, Lapistan'
$ws.Range('B24').Value = 'nyeri pada saluran kencing saat kencing'
$ws.Range('C24').Value = 'Paracetamol, Freedol, Cefadroxil'
$ws.Range('B25').Value = 'kesulitan saat buang air besar'
$ws.Range('B26').Value = 'kulit bernanah akibat infeksi kulit'
$ws.Range('B27').Value = 'nyeri pada tulang dan nyeri pada persendian'
$ws.Range('B28').Value = 'saya nyeri kepala ringan'
$ws.Range('C28').Value = 'Paracetamol, topcilin, Novagesic'

# Match the saved view state: zoom level and active cell selection.
$excel.ActiveWindow.Zoom = 96
$ws.Range('M10').Select()
